$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1082.7142
$ws.Range("I2").Value = 1013.1667
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1013.1667
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -900.1667
$ws.Range("N2").Value = -1726

$ws.Range("H32").Value = 8188.375
$ws.Range("I32").Value = 7289.3335
$ws.Range("J32").Value = 8727.799999999999
$ws.Range("K32").Value = 7289.3335
$ws.Range("L32").Value = 8727.799999999999
$ws.Range("M32").Value = -6963.3335
$ws.Range("N32").Value = -9379.799999999999

$ws.Range("H52").Value = 2000
$ws.Range("I52").Value = 2000
$ws.Range("K52").Value = 6000
$ws.Range("M52").Value = -5840

$ws.Range("H92").Value = 355.7143
$ws.Range("I92").Value = 304.08334
$ws.Range("K92").Value = 304.08334
$ws.Range("M92").Value = 943.91666

$ws.Range("H113").Value = 5666.3335
$ws.Range("I113").Value = 4999.5
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 4999.5
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -1745.5
$ws.Range("N113").Value = -13508

$ws.Range("H135").Value = 801.3043
$ws.Range("I135").Value = 746.8182
$ws.Range("K135").Value = 6721.3638
$ws.Range("M135").Value = -4186.3638

$ws.Range("H137").Value = 3679.8125
$ws.Range("I137").Value = 2398.818
$ws.Range("K137").Value = 7196.454000000001
$ws.Range("M137").Value = -4646.454000000001

$ws.Range("H138").Value = 4227.1577
$ws.Range("I138").Value = 2052.3125
$ws.Range("J138").Value = 5808.864
$ws.Range("K138").Value = 6156.9375
$ws.Range("L138").Value = 17426.592
$ws.Range("M138").Value = -1016.9375
$ws.Range("N138").Value = -27706.592

$ws.Range("H141").Value = 2650.92
$ws.Range("I141").Value = 2289.3333
$ws.Range("J141").Value = 4549.25
$ws.Range("K141").Value = 6867.999899999999
$ws.Range("L141").Value = 13647.75
$ws.Range("M141").Value = -1687.999899999999
$ws.Range("N141").Value = -24007.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2276.5483
$ws.Range("I74").Value = 2326.724
$ws.Range("J74").Value = 1549
$ws.Range("K74").Value = 2326.724
$ws.Range("L74").Value = 1549
$ws.Range("M74").Value = -1452.724
$ws.Range("N74").Value = -3297

$ws.Range("H77").Value = 2276.5483
$ws.Range("I77").Value = 2326.724
$ws.Range("J77").Value = 1549
$ws.Range("K77").Value = 11633.62
$ws.Range("L77").Value = 7745
$ws.Range("M77").Value = -7265.620000000001
$ws.Range("N77").Value = -16481

$ws.Range("H102").Value = 1620.6666
$ws.Range("I102").Value = 1620.6666
$ws.Range("K102").Value = 1620.6666
$ws.Range("M102").Value = 1.333399999999983

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2971.2354
$ws.Range("I134").Value = 2971.2354
$ws.Range("K134").Value = 8913.706200000001
$ws.Range("M134").Value = -6378.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 170.625
$ws.Range("I7").Value = 108
$ws.Range("J7").Value = 275
$ws.Range("K7").Value = 108
$ws.Range("L7").Value = 275
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = -501

$ws.Range("H22").Value = 619.2
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 574
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 574
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -1274

$ws.Range("H31").Value = 3332
$ws.Range("J31").Value = 3665.3333
$ws.Range("L31").Value = 3665.3333
$ws.Range("N31").Value = -4255.3333

$ws.Range("H34").Value = 3332
$ws.Range("J34").Value = 3665.3333
$ws.Range("L34").Value = 3665.3333
$ws.Range("N34").Value = -4069.3333

$ws.Range("H58").Value = 1799.6666
$ws.Range("I58").Value = 1596.5555
$ws.Range("J58").Value = 2409
$ws.Range("K58").Value = 1596.5555
$ws.Range("L58").Value = 2409
$ws.Range("M58").Value = -1393.5555
$ws.Range("N58").Value = -2815

$ws.Range("H86").Value = 35369.89
$ws.Range("J86").Value = 47720.668
$ws.Range("L86").Value = 47720.668
$ws.Range("N86").Value = -49966.668

$ws.Range("H89").Value = 35369.89
$ws.Range("J89").Value = 47720.668
$ws.Range("L89").Value = 238603.34
$ws.Range("N89").Value = -249835.34

$ws.Range("H132").Value = 3795.4
$ws.Range("I132").Value = 3925.7368
$ws.Range("J132").Value = 3382.6667
$ws.Range("K132").Value = 11777.2104
$ws.Range("L132").Value = 10148.0001
$ws.Range("M132").Value = -9247.2104
$ws.Range("N132").Value = -15208.0001

$ws.Range("H134").Value = 3182.4443
$ws.Range("I134").Value = 3182.4443
$ws.Range("K134").Value = 9547.332900000001
$ws.Range("M134").Value = -7012.332900000001

$ws.Range("H136").Value = 1799.6666
$ws.Range("I136").Value = 1596.5555
$ws.Range("J136").Value = 2409
$ws.Range("K136").Value = 4789.666499999999
$ws.Range("L136").Value = 7227
$ws.Range("M136").Value = -2239.666499999999
$ws.Range("N136").Value = -12327

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62.4
$ws.Range("I2").Value = 55
$ws.Range("K2").Value = 330
$ws.Range("M2").Value = -217

$ws.Range("H5").Value = 897.75
$ws.Range("I5").Value = 699.3333
$ws.Range("J5").Value = 1493
$ws.Range("K5").Value = 2097.9999
$ws.Range("L5").Value = 4479
$ws.Range("M5").Value = -1985.9999
$ws.Range("N5").Value = -4703

$ws.Range("H29").Value = 9999
$ws.Range("J29").Value = 9999
$ws.Range("L29").Value = 29997
$ws.Range("N29").Value = -30551

$ws.Range("H34").Value = 1766
$ws.Range("I34").Value = 896.5
$ws.Range("K34").Value = 2689.5
$ws.Range("M34").Value = -2605.5

$ws.Range("H39").Value = 3000
$ws.Range("J39").Value = 3000
$ws.Range("L39").Value = 9000
$ws.Range("N39").Value = -9588

$ws.Range("H40").Value = 152.71428
$ws.Range("J40").Value = 100
$ws.Range("L40").Value = 400
$ws.Range("N40").Value = -538

$ws.Range("H55").Value = 3100
$ws.Range("J55").Value = 3933.3333
$ws.Range("L55").Value = 11799.9999
$ws.Range("N55").Value = -12153.9999

$ws.Range("H118").Value = 4764.3125
$ws.Range("I118").Value = 1229
$ws.Range("K118").Value = 3687
$ws.Range("M118").Value = -2444

$ws.Range("H129").Value = 2900
$ws.Range("I129").Value = 2900
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 8700
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -3700
$ws.Range("N129").ClearContents()

$ws.Range("H135").Value = 897.75
$ws.Range("I135").Value = 699.3333
$ws.Range("J135").Value = 1493
$ws.Range("K135").Value = 6293.9997
$ws.Range("L135").Value = 13437
$ws.Range("M135").Value = -3758.9997
$ws.Range("N135").Value = -18507

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 899
$ws.Range("I122").Value = 899
$ws.Range("K122").Value = 2697
$ws.Range("M122").Value = -247

$ws.Range("H132").Value = 1517.8572
$ws.Range("I132").Value = 806.53845
$ws.Range("K132").Value = 2419.61535
$ws.Range("M132").Value = 110.38465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2425.5
$ws.Range("I93").Value = 2425.5
$ws.Range("K93").Value = 2425.5
$ws.Range("M93").Value = -1177.5

$ws.Range("H122").Value = 2981.6
$ws.Range("I122").Value = 2981.6
$ws.Range("K122").Value = 8944.799999999999
$ws.Range("M122").Value = -6494.799999999999

$ws.Range("H132").Value = 2922.9033
$ws.Range("I132").Value = 1945.0625
$ws.Range("K132").Value = 5835.1875
$ws.Range("M132").Value = -3305.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 443.85715
$ws.Range("I107").Value = 384.83334
$ws.Range("K107").Value = 1154.50002
$ws.Range("M107").Value = 765.4999800000001

$ws.Range("H125").Value = 79997.5
$ws.Range("J125").Value = 79997.5
$ws.Range("L125").Value = 79997.5
$ws.Range("N125").Value = -89837.5

$ws.Range("H132").Value = 2515.2104
$ws.Range("I132").Value = 1758.48
$ws.Range("J132").Value = 3970.4614
$ws.Range("K132").Value = 5275.440000000001
$ws.Range("L132").Value = 11911.3842
$ws.Range("M132").Value = -2745.440000000001
$ws.Range("N132").Value = -16971.3842

$ws.Range("H136").Value = 1123.1875
$ws.Range("I136").Value = 851.5714
$ws.Range("J136").Value = 3024.5
$ws.Range("K136").Value = 2554.7142
$ws.Range("L136").Value = 9073.5
$ws.Range("M136").Value = -4.714200000000346
$ws.Range("N136").Value = -14173.5

Write-Host "Applied scheduled-runner market/profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
